$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '292.42'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.65%'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '31.03'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '0.78%'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '1.53%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07469'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '3.21%'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-5.02%'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.750'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.17%'
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.767'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '1.53%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9207'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2.61%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09304'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '17.10%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1730'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '3.87%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08313'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2.00%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03277'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '5.87%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09942'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.70%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001497'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.01%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.005785'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.82%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.468'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.08%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.148'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '3.55%'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '0.43%'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.01%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.154'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '4.36%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.2124'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-7.84%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04521'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '0.04%'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001219'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '0.68%'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004260'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-3.05%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001300'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-0.22%'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0003394'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '-0.12%'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01616'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '2.53%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04570'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '4.49%'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007461'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '3.01%'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009841'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-2.06%'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1359'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '3.67%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002160'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '6.54%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009284'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-2.33%'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006107'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '6.65%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000750'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-0.10%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.768'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '23.51%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002000'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '-31.09%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '-0.10%'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '-0.10%'
